$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column G
$ws.Range("G1").Value = "Comments"
$ws.Range("G1").HorizontalAlignment = -4108

# Row 4: new DP question
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Count number of binary strings without consecutive 1$([char]0x2019)s"
$ws.Range("C4").Value = "DP"
$ws.Range("D4").Value = "String"
$ws.Range("E4").Value = "easy"
$ws.Range("F4").Value = "GeeksForGeeks"

# Row 5: another new DP question
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Boolean Parenthesization Problem"
$ws.Range("C5").Value = "DP"
$ws.Range("D5").Value = "String"
$ws.Range("E5").Value = "medium"
$ws.Range("F5").Value = "GeeksForGeeks"

# Column width adjustments (ColumnWidth is pixel-quantized by the engine;
# these inputs land on the closest achievable stored width)
$ws.Columns.Item(2).ColumnWidth = 53.1666666666667
$ws.Columns.Item(7).ColumnWidth = 22.8333333333333

# Update selection to match target state
$ws.Range("F8").Select()
